$wb = $excel.ActiveWorkbook

# --- Add the new sheet "Max_land_usage_global" after "Global_Min_RES_elec_penetration" ---
$afterSheet = $wb.Worksheets.Item("Global_Min_RES_elec_penetration")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "Max_land_usage_global"

# --- Header row (row 1): Years + technology names ---
$headers = @("Years","Natural_gas_supply","Oil_supply","BW_supply","Geo_PP","PV_PP","Wind_PP","Hydro_PP","HFO_PP","OCGT_PP","BW_PP","Elec_transmission_distribution")
for ($col = 1; $col -le 12; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- Column A (rows 2-12): Y0..Y10 ---
$years = @("Y0","Y1","Y2","Y3","Y4","Y5","Y6","Y7","Y8","Y9","Y10")
for ($row = 2; $row -le 12; $row++) {
    $newSheet.Cells.Item($row, 1).Value = $years[$row - 2]
}

# --- Data cells B2:L12 = 1E+30 ---
$bigValue = [double]"1E+30"
for ($row = 2; $row -le 12; $row++) {
    for ($col = 2; $col -le 12; $col++) {
        $newSheet.Cells.Item($row, $col).Value = $bigValue
    }
}

# --- Formatting: bold, bordered, centered header row + column A ---
$headerRange = $newSheet.Range("A1:L1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$colARange = $newSheet.Range("A2:A12")
$colARange.Font.Bold = $true
$colARange.Borders.LineStyle = 1
$colARange.HorizontalAlignment = -4108
$colARange.VerticalAlignment = -4160

# --- Sheet view: selection on B2:L12, tab selected ---
$newSheet.Range("B2:L12").Select()

$wb.Worksheets.Item("Max_newcap_global").Select()
$newSheet.Select()
